# #5: property boat&car done
# Rebuild the "汽車" (car) sheet (sheet3) with the full standard property
# column layout (name, capacity, owner, register_date, register_reason,
# acquire_value, property_category, category, date, legislator_name,
# legislator_id, source_file, index) instead of the old truncated layout.

$wb  = $excel.ActiveWorkbook
$ws  = $wb.Worksheets.Item("汽車")

# ---- Header row (row 1) --------------------------------------------------
$ws.Range("B1").Value2 = "name"
$ws.Range("C1").Value2 = "capacity"
$ws.Range("D1").Value2 = "owner"
$ws.Range("E1").Value2 = "register_date"
$ws.Range("F1").Value2 = "register_reason"
$ws.Range("G1").Value2 = "acquire_value"
$ws.Range("H1").Value2 = "property_category"
$ws.Range("I1").Value2 = "category"
$ws.Range("J1").Value2 = "date"
$ws.Range("K1").Value2 = "legislator_name"
$ws.Range("L1").Value2 = "legislator_id"
$ws.Range("M1").Value2 = "source_file"
$ws.Range("N1").Value2 = "index"

# Match the formatting already used across row 1 (bold / bordered / centered
# header style) for the newly added H1:N1 header cells.
$ws.Range("B1").Copy()
$ws.Range("H1:N1").PasteSpecial(-4122)

# ---- Row 2 (Livina) -------------------------------------------------------
$ws.Range("B2").Value2 = "日產Livina"
$ws.Range("E2").Value2 = "96年12月"
$ws.Range("H2").Value2 = "land"
$ws.Range("I2").Value2 = "normal"
# Leading apostrophe forces text (otherwise "2012-04-30" is auto-parsed as a date).
$ws.Range("J2").Value2 = "'2012-04-30"
$ws.Range("K2").Value2 = "邱文彥"
$ws.Range("L2").Value2 = 1743
$ws.Range("M2").Value2 = "tmpa1171"
$ws.Range("N2").Value2 = 42

# ---- Row 3 (Sentra) --------------------------------------------------------
$ws.Range("B3").Value2 = "曰產Sentra"
$ws.Range("E3").Value2 = "91年8月"
$ws.Range("H3").Value2 = "land"
$ws.Range("I3").Value2 = "normal"
$ws.Range("J3").Value2 = "'2012-04-30"
$ws.Range("K3").Value2 = "邱文彥"
$ws.Range("L3").Value2 = 1743
$ws.Range("M3").Value2 = "tmpa1171"
$ws.Range("N3").Value2 = 43
